$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "만드려면" -> "만들려면" in the pizza question (row 18, column B)
$ws.Range("B18").Value = "bm² 넓이의 피자를 만들려면 치즈 akg가 필요합니다. 1m² 넓이의 피자를 만들기 위해 필요한 치즈는 몇 kg인지 구해 보세요."

# Fix typo "만드려면" -> "만들려면" in the kimchi question (row 19, column B)
$ws.Range("B19").Value = "김치 bkg를 만들려면 고춧가루 akg가 필요합니다. 김치 1kg를 만들기 위해 필요한 고춧가루는 몇 kg인지 구해 보세요."

# Update the active selection to B24 (matches sheetView selection change in the diff)
$ws.Range("B24").Select()
